# Updated Results to include averages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H width (matches the 14.75-character width added to the sheet)
$ws.Columns.Item(8).ColumnWidth = 14

# --- Row 1: header label for the new "Project Accuracy" column ---
$ws.Range("F1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Project Accuracy:"

# --- Row 2: "Top 3 Correct" label + value ---
$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H2").Value = "Top 3 Correct:"
$ws.Range("F2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").Value = 0.46

# --- Row 3: "Total Predictions" label + value ---
$ws.Range("F2").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").Value = "Total Predictions:"
$ws.Range("F2").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value = 0.65

# --- Row 78: project averages row ---
$ws.Range("F2").Copy()
$ws.Range("D78").PasteSpecial(-4122)
$ws.Range("D78").Value = 0.46
$ws.Range("F2").Copy()
$ws.Range("E78").PasteSpecial(-4122)
$ws.Range("E78").Value = 0.65
$ws.Range("F2").Copy()
$ws.Range("G78").PasteSpecial(-4122)
$ws.Range("G78").Value = "Averages"
